$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 49, pushing the existing rows 49:67 down to 50:68.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly record.
$ws.Range("A49").Value = 8
$ws.Range("B49").Value = "Terminal La Palmera de La Serena"
$ws.Range("C49").Value = "Coquimbo"
$ws.Range("D49").Value = 44839
$ws.Range("E49").Value = 4
$ws.Range("F49").Value = 100114007
$ws.Range("G49").Value = "Jengibre"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 520
$ws.Range("K49").Value = 14000
$ws.Range("L49").Value = 15000
$ws.Range("M49").Value = 14500
$ws.Range("N49").Value = "$/caja 13 kilos"
$ws.Range("O49").Value = "Perú"
$ws.Range("P49").Value = 1115
$ws.Range("Q49").Value = 13
$ws.Range("R49").Value = "Hortaliza"
